# The password test value on the "LoginTestData" sheet is updated from
# "Faris198407" to "admin123", and the workbook's active sheet/selection
# moves from "New_Employee" (B2/E20) to "LoginTestData" cell B2.

$wb = $excel.ActiveWorkbook

$loginSheet = $wb.Worksheets.Item("LoginTestData")
$loginSheet.Range("B2").Value = "admin123"

$loginSheet.Activate()
$loginSheet.Range("B2").Select()
